{"js": "// Add group member names and contribution levels to the empty cells of\n// the \"Group Contribution Form\" table (rows 2-4, i.e. the first three\n// data rows under the header row).\n//\n// Each target cell's first paragraph is currently empty (it only carries\n// paragraph-level formatting in <w:pPr><w:rPr>). We insert a run of text\n// into that existing paragraph (preserving the paragraph itself) and set\n// the run's font to match the paragraph's own font so the inserted text\n// looks the same as it would if typed directly into the document.\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// [rowIndex, columnIndex, paragraphIndex, text, fontName] -- all 0-based\n// except fontName.\nconst entries = [\n  [1, 0, 0, \"Vanisha Oree\", \"Calibri\"],\n  [1, 1, 0, \"+1\", \"Calibri\"],\n  [2, 0, 0, \"David Hood\", \"Calibri\"],\n  [2, 1, 0, \"+1\", \"Calibri\"],\n  [2, 2, 0, \"David Hood\", \"Informal Roman\"],\n  [3, 0, 0, \"Loving-Grace Mawire\", \"Calibri\"],\n  [3, 1, 0, \"-2\", \"Calibri\"],\n];\n\nfor (const [rowIndex, colIndex, paraIndex, text, fontName] of entries) {\n  const cell = table.getCell(rowIndex, colIndex);\n  const paragraphs = cell.body.paragraphs;\n  paragraphs.load(\"items\");\n  await context.sync();\n\n  const paragraph = paragraphs.items[paraIndex];\n  const range = paragraph.insertText(text, Word.InsertLocation.end);\n  range.font.name = fontName;\n  await context.sync();\n}\n", "ps1": "# Add group member names and contribution levels to the empty cells of\n# the \"Group Contribution Form\" table (the first three data rows right\n# below the header row).\n#\n# Each target cell's first paragraph is currently empty (it only carries\n# paragraph-level formatting). We insert the text into that existing\n# paragraph (leaving the paragraph itself intact) and then set the\n# inserted run's font to match the paragraph's own font.\n\n$d = $word.ActiveDocument\n$table = $d.Tables.Item(1)\n\n# row (1-based), column (1-based), paragraph-in-cell (1-based), text, font name\n$entries = @(\n    @(2, 1, 1, \"Vanisha Oree\", \"Calibri\"),\n    @(2, 2, 1, \"+1\", \"Calibri\"),\n    @(3, 1, 1, \"David Hood\", \"Calibri\"),\n    @(3, 2, 1, \"+1\", \"Calibri\"),\n    @(3, 3, 1, \"David Hood\", \"Informal Roman\"),\n    @(4, 1, 1, \"Loving-Grace Mawire\", \"Calibri\"),\n    @(4, 2, 1, \"-2\", \"Calibri\")\n)\n\nforeach ($entry in $entries) {\n    $rowIndex = $entry[0]\n    $colIndex = $entry[1]\n    $paraIndex = $entry[2]\n    $text = $entry[3]\n    $fontName = $entry[4]\n\n    # Re-fetch the table/cell/paragraph fresh each time so the handle is\n    # never stale after the previous iteration's structural edit.\n    $table = $d.Tables.Item(1)\n    $cell = $table.Cell($rowIndex, $colIndex)\n    $para = $cell.Range.Paragraphs.Item($paraIndex)\n    $rng = $para.Range\n    $rng.InsertBefore($text)\n\n    # Re-fetch again before touching formatting, since InsertBefore just\n    # mutated the document and invalidated this range/paragraph handle.\n    $table = $d.Tables.Item(1)\n    $cell = $table.Cell($rowIndex, $colIndex)\n    $para = $cell.Range.Paragraphs.Item($paraIndex)\n    $para.Range.Font.Name = $fontName\n}\n"}
